# This edit performs a cyclic rotation of the species-observation data
# stored in rows 2-4 of the active sheet:
#   new row 2 = old row 3
#   new row 3 = old row 4
#   new row 4 = old row 2
# Only the columns that actually vary per-row are rotated: A, B, E, F, G, H,
# Q, R, AO. (All the other columns are identical across rows 2-4 already, so
# rotating them would be a no-op, but we leave them untouched regardless.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AO")

# Snapshot current values for rows 2, 3, 4 before writing anything, so the
# rotation doesn't clobber source data it still needs to read.
$row2 = @{}
$row3 = @{}
$row4 = @{}

foreach ($col in $cols) {
    $row2[$col] = $ws.Range("$col" + "2").Value2
    $row3[$col] = $ws.Range("$col" + "3").Value2
    $row4[$col] = $ws.Range("$col" + "4").Value2
}

# Apply the rotation: row2 <- row3, row3 <- row4, row4 <- row2(original)
foreach ($col in $cols) {
    $ws.Range("$col" + "2").Value2 = $row3[$col]
    $ws.Range("$col" + "3").Value2 = $row4[$col]
    $ws.Range("$col" + "4").Value2 = $row2[$col]
}
